# Update the localization status text from "Ready for handoff" to "In Translation"
# on every sheet that reports per-language status, and shrink the now-narrower
# status columns to match (report regenerated for archive).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New column width, expressed in COM "characters" units, that results in the
# narrower OOXML column width used by the regenerated report.
$newColumnWidth = 12.5

# Overview sheet: zh-cn (col E) and de-de (col F) status cells + their columns.
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: Status column (col C).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: Status column (col C).
$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
